$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (will be re-added on their new, shifted cells)
$ws.Hyperlinks.Delete()

# Delete column F ("jezyk") - cells to the right shift left by one column
$ws.Columns.Item(6).Delete()

# Re-create hyperlinks on their new locations (shifted from I to H)
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:olo@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "olo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ala@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ala@gmail.com")

# Update selection to match the post-edit state
$ws.Range("F1").Select()
